$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from the existing "Unnamed: 28" header cell (AC1)
# so the new headers match the existing bold/border/centered formatting.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the Wins/Losses/Ties values for every data row (2-43)
$ws.Range("AD2:AD43").Value = 84
$ws.Range("AE2:AE43").Value = 78
$ws.Range("AF2:AF43").Value = 0
